# Update "Favorites / views" (column F) counts across the workbook's sheets.
# Mirrors a refreshed data pull: most counters ticked up slightly.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 243
$ws1.Range("F4").Value  = 830
$ws1.Range("F5").Value  = 242
$ws1.Range("F6").Value  = 411
$ws1.Range("F7").Value  = 585
$ws1.Range("F8").Value  = 218
$ws1.Range("F10").Value = 340
$ws1.Range("F11").Value = 143
$ws1.Range("F12").Value = 642
$ws1.Range("F13").Value = 85
$ws1.Range("F14").Value = 1795
$ws1.Range("F15").Value = 344
$ws1.Range("F16").Value = 2943
$ws1.Range("F17").Value = 314
$ws1.Range("F19").Value = 50

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value  = 476
$ws2.Range("F13").Value = 90
$ws2.Range("F14").Value = 39

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5313
$ws3.Range("F3").Value = 316
$ws3.Range("F4").Value = 243

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5313
$ws4.Range("F4").Value  = 316
$ws4.Range("F6").Value  = 243
$ws4.Range("F7").Value  = 243
$ws4.Range("F12").Value = 476
$ws4.Range("F13").Value = 830
$ws4.Range("F16").Value = 242
$ws4.Range("F17").Value = 411
$ws4.Range("F18").Value = 585
$ws4.Range("F19").Value = 218
$ws4.Range("F22").Value = 340
$ws4.Range("F23").Value = 143
$ws4.Range("F26").Value = 642
$ws4.Range("F27").Value = 85
$ws4.Range("F28").Value = 90
$ws4.Range("F29").Value = 1795
$ws4.Range("F30").Value = 344
$ws4.Range("F31").Value = 2944
$ws4.Range("F32").Value = 39
$ws4.Range("F33").Value = 314
$ws4.Range("F35").Value = 50
